{"js": "// Fix inaccuracies in the library-description bullet list:\n//  1) \"... \u0441 \u0448\u0438\u0444\u0440\u043e\u043c Base64\"         -> \"... \u0441 \u043a\u043e\u0434\u0438\u0440\u043e\u0432\u043a\u043e\u0439 Base64\"\n//  2) \"... \u0441 \u0441\u0438\u043c\u043c\u0435\u0442\u0440\u0438\u0447\u043d\u044b\u043c\u0438 \u0448\u0438\u0444\u0440\u0430\u043c\u0438\" -> \"... \u0441 \u0441\u0438\u043c\u043c\u0435\u0442\u0440\u0438\u0447\u043d\u044b\u043c \u0448\u0438\u0444\u0440\u043e\u0432\u0430\u043d\u0438\u0435\u043c\"\n//  3) \"... \u0441 \u041f\u041e Gpg4win\"            -> \"... \u0441 GPG\"\n\nconst body = context.document.body;\n\n// --- Fix 1: \"base64 ... \u0441 \u0448\u0438\u0444\u0440\u043e\u043c Base64\" -> \"... \u0441 \u043a\u043e\u0434\u0438\u0440\u043e\u0432\u043a\u043e\u0439 Base64\" ---\nconst fix1 = body.search(\" \u0441 \u0448\u0438\u0444\u0440\u043e\u043c \", { matchCase: true, matchWholeWord: false });\nfix1.load(\"items\");\n\n// --- Fix 2: \"pycryptodome ... \u0441 \u0441\u0438\u043c\u043c\u0435\u0442\u0440\u0438\u0447\u043d\u044b\u043c\u0438 \u0448\u0438\u0444\u0440\u0430\u043c\u0438\" -> \"... \u0441 \u0441\u0438\u043c\u043c\u0435\u0442\u0440\u0438\u0447\u043d\u044b\u043c \u0448\u0438\u0444\u0440\u043e\u0432\u0430\u043d\u0438\u0435\u043c\" ---\nconst fix2 = body.search(\" \u0434\u043b\u044f \u0440\u0430\u0431\u043e\u0442\u044b \u0441 \u0441\u0438\u043c\u043c\u0435\u0442\u0440\u0438\u0447\u043d\u044b\u043c\u0438 \u0448\u0438\u0444\u0440\u0430\u043c\u0438\", { matchCase: true, matchWholeWord: false });\nfix2.load(\"items\");\n\nawait context.sync();\n\nif (fix1.items.length > 0) {\n  fix1.items[0].insertText(\" \u0441 \u043a\u043e\u0434\u0438\u0440\u043e\u0432\u043a\u043e\u0439 \", Word.InsertLocation.replace);\n}\nif (fix2.items.length > 0) {\n  fix2.items[0].insertText(\" \u0434\u043b\u044f \u0440\u0430\u0431\u043e\u0442\u044b \u0441 \u0441\u0438\u043c\u043c\u0435\u0442\u0440\u0438\u0447\u043d\u044b\u043c \u0448\u0438\u0444\u0440\u043e\u0432\u0430\u043d\u0438\u0435\u043c\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n\n// --- Fix 3: \"python-gnupg ... \u0434\u043b\u044f \u0440\u0430\u0431\u043e\u0442\u044b \u0441 \u041f\u041e Gpg4win\" -> \"... \u0434\u043b\u044f \u0440\u0430\u0431\u043e\u0442\u044b \u0441 GPG\" ---\n// Done in two narrow steps so the unrelated \"python-gnupg\" spell-check markers\n// earlier in the same sentence are left untouched.\nconst poFix = body.search(\" \u0434\u043b\u044f \u0440\u0430\u0431\u043e\u0442\u044b \u0441 \u041f\u041e \", { matchCase: true, matchWholeWord: false });\npoFix.load(\"items\");\nawait context.sync();\n\nif (poFix.items.length > 0) {\n  poFix.items[0].insertText(\" \u0434\u043b\u044f \u0440\u0430\u0431\u043e\u0442\u044b \u0441 \", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Re-locate the paragraph via the still-unique \"gnupg\" anchor, then replace\n// \"Gpg4win\" only inside that paragraph so other \"Gpg4win\" mentions elsewhere\n// in the document are left alone.\nconst anchor = body.search(\"gnupg\", { matchCase: true, matchWholeWord: false });\nanchor.load(\"items\");\nawait context.sync();\n\nif (anchor.items.length > 0) {\n  const targetParagraph = anchor.items[0].paragraphs.getFirst();\n  const gpgFix = targetParagraph.search(\"Gpg4win\", { matchCase: true, matchWholeWord: false });\n  gpgFix.load(\"items\");\n  await context.sync();\n\n  if (gpgFix.items.length > 0) {\n    gpgFix.items[0].insertText(\"GPG\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fix inaccuracies in the library-description bullet list:\n#  1) \"... \u0441 \u0448\u0438\u0444\u0440\u043e\u043c Base64\"         -> \"... \u0441 \u043a\u043e\u0434\u0438\u0440\u043e\u0432\u043a\u043e\u0439 Base64\"\n#  2) \"... \u0441 \u0441\u0438\u043c\u043c\u0435\u0442\u0440\u0438\u0447\u043d\u044b\u043c\u0438 \u0448\u0438\u0444\u0440\u0430\u043c\u0438\" -> \"... \u0441 \u0441\u0438\u043c\u043c\u0435\u0442\u0440\u0438\u0447\u043d\u044b\u043c \u0448\u0438\u0444\u0440\u043e\u0432\u0430\u043d\u0438\u0435\u043c\"\n#  3) \"... \u0441 \u041f\u041e Gpg4win\"            -> \"... \u0441 GPG\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-DocText \" \u0441 \u0448\u0438\u0444\u0440\u043e\u043c \" \" \u0441 \u043a\u043e\u0434\u0438\u0440\u043e\u0432\u043a\u043e\u0439 \"\nReplace-DocText \" \u0434\u043b\u044f \u0440\u0430\u0431\u043e\u0442\u044b \u0441 \u0441\u0438\u043c\u043c\u0435\u0442\u0440\u0438\u0447\u043d\u044b\u043c\u0438 \u0448\u0438\u0444\u0440\u0430\u043c\u0438\" \" \u0434\u043b\u044f \u0440\u0430\u0431\u043e\u0442\u044b \u0441 \u0441\u0438\u043c\u043c\u0435\u0442\u0440\u0438\u0447\u043d\u044b\u043c \u0448\u0438\u0444\u0440\u043e\u0432\u0430\u043d\u0438\u0435\u043c\"\nReplace-DocText \" \u0434\u043b\u044f \u0440\u0430\u0431\u043e\u0442\u044b \u0441 \u041f\u041e Gpg4win\" \" \u0434\u043b\u044f \u0440\u0430\u0431\u043e\u0442\u044b \u0441 GPG\"\n"}
